$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1802120141342756
$ws.Range("C2").Value = 0.5618374558303887
$ws.Range("J2").Value = 0.03886925795053003
$ws.Range("P2").Value = 0.1342756183745583
$ws.Range("S2").Value = 0.08480565371024736
$ws.Range("B3").Value = 0.006134969325153374
$ws.Range("C3").Value = 0.03067484662576687
$ws.Range("J3").Value = 0.049079754601227
$ws.Range("P3").Value = 0.7361963190184049
$ws.Range("S3").Value = 0.1779141104294479
$ws.Range("J4").Value = 0.1025641025641026
$ws.Range("O4").Value = 0.02564102564102564
$ws.Range("P4").Value = 0.6410256410256411
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("B6").Value = 0.07547169811320754
$ws.Range("D6").Value = 0.009433962264150943
$ws.Range("F6").Value = 0.04716981132075472
$ws.Range("J6").Value = 0.3207547169811321
$ws.Range("O6").Value = 0.01886792452830189
$ws.Range("Q6").Value = 0.1650943396226415
$ws.Range("R6").Value = 0.08018867924528301
$ws.Range("S6").Value = 0.2830188679245283
$ws.Range("B7").Value = 0.07926829268292683
$ws.Range("D7").Value = 0.02439024390243903
$ws.Range("E7").Value = 0.006097560975609756
$ws.Range("F7").Value = 0.06097560975609756
$ws.Range("J7").Value = 0.1280487804878049
$ws.Range("O7").Value = 0.02439024390243903
$ws.Range("Q7").Value = 0.1768292682926829
$ws.Range("R7").Value = 0.07317073170731707
$ws.Range("S7").Value = 0.4268292682926829
$ws.Range("B8").Value = 0.08040201005025126
$ws.Range("D8").Value = 0.01005025125628141
$ws.Range("E8").Value = 0.002512562814070352
$ws.Range("F8").Value = 0.06532663316582915
$ws.Range("J8").Value = 0.1055276381909548
$ws.Range("O8").Value = 0.02010050251256281
$ws.Range("Q8").Value = 0.1834170854271357
$ws.Range("R8").Value = 0.09798994974874371
$ws.Range("S8").Value = 0.4346733668341708
$ws.Range("B9").Value = 0.1155555555555556
$ws.Range("D9").Value = 0.02666666666666667
$ws.Range("F9").Value = 0.04444444444444445
$ws.Range("J9").Value = 0.1555555555555556
$ws.Range("O9").Value = 0.01333333333333333
$ws.Range("Q9").Value = 0.1466666666666667
$ws.Range("R9").Value = 0.12
$ws.Range("S9").Value = 0.3777777777777778
$ws.Range("B10").Value = 0.10546875
$ws.Range("D10").Value = 0.0203125
$ws.Range("F10").Value = 0.06484375000000001
$ws.Range("J10").Value = 0.15390625
$ws.Range("O10").Value = 0.0109375
$ws.Range("Q10").Value = 0.2078125
$ws.Range("R10").Value = 0.09531249999999999
$ws.Range("S10").Value = 0.34140625
$ws.Range("G11").Value = 0.1504065040650407
$ws.Range("J11").Value = 0.0975609756097561
$ws.Range("K11").Value = 0.2032520325203252
$ws.Range("L11").Value = 0.5365853658536586
$ws.Range("S11").Value = 0.01219512195121951
$ws.Range("G12").Value = 0.8041958041958042
$ws.Range("J12").Value = 0.0979020979020979
$ws.Range("K12").Value = 0.02097902097902098
$ws.Range("L12").Value = 0.05594405594405594
$ws.Range("S12").Value = 0.02097902097902098
$ws.Range("G13").Value = 0.6774193548387096
$ws.Range("J13").Value = 0.3225806451612903
$ws.Range("F15").Value = 0.01724137931034483
$ws.Range("H15").Value = 0.1724137931034483
$ws.Range("I15").Value = 0.06896551724137931
$ws.Range("J15").Value = 0.3793103448275862
$ws.Range("K15").Value = 0.04022988505747126
$ws.Range("M15").Value = 0.01724137931034483
$ws.Range("O15").Value = 0.06321839080459771
$ws.Range("S15").Value = 0.2413793103448276
$ws.Range("F16").Value = 0.03409090909090909
$ws.Range("H16").Value = 0.1477272727272727
$ws.Range("I16").Value = 0.1306818181818182
$ws.Range("J16").Value = 0.4715909090909091
$ws.Range("K16").Value = 0.08522727272727272
$ws.Range("M16").Value = 0.01136363636363636
$ws.Range("O16").Value = 0.06818181818181818
$ws.Range("S16").Value = 0.05113636363636364
$ws.Range("F17").Value = 0.02540415704387991
$ws.Range("H17").Value = 0.187066974595843
$ws.Range("I17").Value = 0.115473441108545
$ws.Range("J17").Value = 0.4018475750577367
$ws.Range("K17").Value = 0.09930715935334873
$ws.Range("M17").Value = 0.01616628175519631
$ws.Range("O17").Value = 0.05311778290993072
$ws.Range("S17").Value = 0.1016166281755196
$ws.Range("F18").Value = 0.04608294930875576
$ws.Range("H18").Value = 0.1751152073732719
$ws.Range("I18").Value = 0.1290322580645161
$ws.Range("J18").Value = 0.4377880184331797
$ws.Range("K18").Value = 0.06912442396313365
$ws.Range("M18").Value = 0.02304147465437788
$ws.Range("O18").Value = 0.04147465437788019
$ws.Range("S18").Value = 0.07834101382488479
$ws.Range("F19").Value = 0.0168697282099344
$ws.Range("H19").Value = 0.2071227741330834
$ws.Range("I19").Value = 0.1021555763823805
$ws.Range("J19").Value = 0.4236176194939081
$ws.Range("K19").Value = 0.1030927835051546
$ws.Range("M19").Value = 0.014058106841612
$ws.Range("O19").Value = 0.05435801312089972
$ws.Range("S19").Value = 0.07872539831302718
